$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded. Insert a row at 187 (pushing
# the existing rows 187:225 down to 188:226) and populate it with the new
# data point, matching the layout of the surrounding rows.
$ws.Rows.Item(187).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

$ws.Range("A187").Value = 6
$ws.Range("B187").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C187").Value = "Metropolitana"
$ws.Range("D187").Value = 44785
$ws.Range("E187").Value = 13
$ws.Range("F187").Value = 100112029
$ws.Range("G187").Value = "Orégano"
$ws.Range("H187").Value = "Sin especificar"
$ws.Range("I187").Value = "Primera"
$ws.Range("J187").Value = 53
$ws.Range("K187").Value = 15000
$ws.Range("L187").Value = 16000
$ws.Range("M187").Value = 15528
$ws.Range("N187").Value = '$/docena de atados'
$ws.Range("O187").Value = "Región Metropolitana"
$ws.Range("P187").Value = 5176
$ws.Range("Q187").Value = 3
$ws.Range("R187").Value = "Hortaliza"

# Carry the date-number formatting from the row above onto the new row's
# date cell so it keeps rendering as a date rather than a raw serial.
$ws.Range("D186").Copy()
$ws.Range("D187").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
